# Append 16 new case-disposition rows (21TRD09200 / Bunner) to Sheet1,
# rows 1395-1410, columns A-K (J/K populated only on some rows), per the
# commit "Updated tests for fine only."

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Helper: decide whether a text value would be auto-coerced by Excel into
# a number (e.g. "4510.111", "$ 0") and needs a leading apostrophe to stay
# text, matching how the rest of the sheet stores everything as text.
function Needs-TextForce($v) {
    if ($v -match '^\s*[+-]?(\d+\.?\d*|\.\d+)\s*$') { return $true }
    if ($v -match '^\s*\$\s*[\d,]+\.?\d*\s*$') { return $true }
    return $false
}

function Set-TextCell($sheet, $targetRow, $targetCol, $value) {
    if ($null -eq $value) { return }
    if (Needs-TextForce $value) {
        $sheet.Cells.Item($targetRow, $targetCol).Value2 = "'" + $value
    } else {
        $sheet.Cells.Item($targetRow, $targetCol).Value2 = $value
    }
}

# Each entry: A,B,C,D,E,F,G,H,I,J,K ($null means the cell is left empty/absent)
# (built with ArrayList.Add so each row stays a distinct nested array instead
# of being flattened into the outer collection)
$rows = New-Object System.Collections.ArrayList
$rows.Add(@("21TRD09200","Bunner","DUS Ucm","4510.111","UCM","Guilty","Guilty","$ 0","$ 0",$null,$null)) | Out-Null
$rows.Add(@("21TRD09200","Bunner","Operating W/o A Valid Ol - Ucm","4510.12","UCM","Guilty","Guilty","$ 0","$ 0",$null,$null)) | Out-Null
$rows.Add(@("21TRD09200","Bunner","Failure To Reinstate License Ucm 1-2/3yrs","4510.21A*","UCM","Guilty","Guilty","$ 0","$ 0",$null,$null)) | Out-Null
$rows.Add(@("21TRD09200","Bunner","Failure To File Registration","4503.11","MM","Guilty","Guilty","$ 0","$ 0",$null,$null)) | Out-Null
$rows.Add(@("21TRD09200","Bunner","DUS Ucm","4510.111","UCM","Guilty","Guilty","$ 0","$ 0","None","None")) | Out-Null
$rows.Add(@("21TRD09200","Bunner","Operating W/o A Valid Ol - Ucm","4510.12","UCM","Guilty","Guilty","$ 0","$ 0","None","None")) | Out-Null
$rows.Add(@("21TRD09200","Bunner","Failure To Reinstate License Ucm 1-2/3yrs","4510.21A*","UCM","Guilty","Guilty","$ 0","$ 0","None","None")) | Out-Null
$rows.Add(@("21TRD09200","Bunner","Failure To File Registration","4503.11","MM","Guilty","Guilty","$ 0","$ 0","None","None")) | Out-Null
$rows.Add(@("21TRD09200","Bunner","DUS Ucm","4510.111","UCM","No Contest","Guilty","$ 0","$ 0",$null,$null)) | Out-Null
$rows.Add(@("21TRD09200","Bunner","Operating W/o A Valid Ol - Ucm","4510.12","UCM","No Contest","Guilty","$ 0","$ 0",$null,$null)) | Out-Null
$rows.Add(@("21TRD09200","Bunner","Failure To Reinstate License Ucm 1-2/3yrs","4510.21A*","UCM","No Contest","Guilty","$ 0","$ 0",$null,$null)) | Out-Null
$rows.Add(@("21TRD09200","Bunner","Failure To File Registration","4503.11","MM","No Contest","Guilty","$ 0","$ 0",$null,$null)) | Out-Null
$rows.Add(@("21TRD09200","Bunner","DUS Ucm","4510.111","UCM","No Contest","Guilty","$ 0","$ 0","None","None")) | Out-Null
$rows.Add(@("21TRD09200","Bunner","Operating W/o A Valid Ol - Ucm","4510.12","UCM","No Contest","Guilty","$ 0","$ 0","None","None")) | Out-Null
$rows.Add(@("21TRD09200","Bunner","Failure To Reinstate License Ucm 1-2/3yrs","4510.21A*","UCM","No Contest","Guilty","$ 0","$ 0","None","None")) | Out-Null
$rows.Add(@("21TRD09200","Bunner","Failure To File Registration","4503.11","MM","No Contest","Guilty","$ 0","$ 0","None","None")) | Out-Null

$startRow = 1395
$r = $startRow
foreach ($dataRow in $rows) {
    for ($c = 0; $c -lt $dataRow.Length; $c++) {
        $colNum = $c + 1
        $cellVal = $dataRow[$c]
        Set-TextCell $ws $r $colNum $cellVal
    }
    $r = $r + 1
}

$lastRow = $r - 1
Write-Output "Added $($rows.Count) rows starting at $startRow; new dimension should be A2:K$lastRow"
